# Weekly update: insert a new price-report row for Haba at "Vega Modelo de
# Temuco" as row 78, pushing the existing rows 78-92 down to 79-93.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 78 (shifts 78..92 down to 79..93,
# and copies formatting - e.g. the date style on column D - from the row above).
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new record's data.
$ws.Cells.Item(78, 1).Value = 10
$ws.Cells.Item(78, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(78, 3).Value = "La Araucanía"
$ws.Cells.Item(78, 4).Value = 45204
$ws.Cells.Item(78, 5).Value = 9
$ws.Cells.Item(78, 6).Value = 100112026
$ws.Cells.Item(78, 7).Value = "Haba"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 45
$ws.Cells.Item(78, 11).Value = 15000
$ws.Cells.Item(78, 12).Value = 15000
$ws.Cells.Item(78, 13).Value = 15000
$ws.Cells.Item(78, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(78, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(78, 16).Value = 600
$ws.Cells.Item(78, 17).Value = 25
$ws.Cells.Item(78, 18).Value = "Hortaliza"
